$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row coloring (per commit message "added colors to rows") ---
# Row 7 and Row 12 -> blue (FF29A3CC)
$ws.Range("A7:J7").Interior.Color = 13411113
$ws.Range("A12:J12").Interior.Color = 13411113

# Row 13 -> orange/yellow (FFFFCC66)
$ws.Range("A13:J13").Interior.Color = 6737151

# Row 14 -> red (FFDF5E5E)
$ws.Range("A14:J14").Interior.Color = 6184671

# I14 changes from 0 to 1
$ws.Range("I14").Value = 1

# --- B19: shared-string " " -> boolean FALSE ---
# B19 sits inside the merged range A19:G19, so it must be unmerged to
# become independently addressable, written to, then the merge restored.
$ws.Range("A19:G19").UnMerge()
$ws.Range("B19").Value = $false

# --- Formula fixes: FLOOR(x,1,1) -> FLOOR(x,1) (drop stray 3rd arg) ---
$ws.Range("B22").Formula = "=FLOOR(F17/8,1)&"".""&FLOOR(MOD(F17,8),1)&"".""&(MOD(F17,8)-FLOOR(MOD(F17,8),1))*60"
$ws.Range("B23").Formula = "=FLOOR(H19,1)&"".""&(H19-FLOOR(H19,1))*8&"".0"""
$ws.Range("B24").Formula = "=FLOOR(I19,1)&"".""&(I19-FLOOR(I19,1))*8&"".0"""
$ws.Range("B27").Formula = "=FLOOR(K27/8,1)&"".""&FLOOR(MOD(K27,8),1)&"".""&(MOD(K27,8)-FLOOR(MOD(K27,8),1))*60"
